$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Custom Statuses")

$newRows = @(
    @("Attempting to softlock irl by biting a cactus", "depressedmonke_"),
    @("Mono gang", "waltuh"),
    @("Playing Celestial Custodian in Graveyard to correct its turn order", "depressedmonke_"),
    @("Ten likes and Huyn draws Swabbie with huge muscles and abs", "lumpymilktea"),
    @("Taking deckbuilding lessons from Crimpton 🛹", "depressedmonke_"),
    @("Ferb, I know what we're gonna do today!", "Tbone")
)

$startRow = 285
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

$wb.Save()
